# [#3] When cell contains a formula, value is always returned.
#
# Adds a new column F ("E") whose rows compute A+C for each data row,
# so a formula-bearing column exists in the test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "E"

# New formula column: F = A + C for each data row
$ws.Range("F2").Formula = "=A2+C2"
$ws.Range("F3").Formula = "=A3+C3"
$ws.Range("F4").Formula = "=A4+C4"

# Match the saved selection state from the target workbook
$ws.Range("I11").Select() | Out-Null
